$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.275.45"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.577.02"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.25"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.799.35"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "1.578.38"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "27.260.18"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.80"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "0.0₃0687"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0465"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.411.57"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.941"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.517"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.88"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "1.711.10"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "0.0₇0992"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.11%  "
